$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 formatting: the merged note cell's row got taller (33pt -> 41.25pt) ---
$ws.Rows.Item(42).RowHeight = 41.25

# --- View state: scroll the window up a bit and select the full row 42 ---
$ws.Range("A42:XFD42").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
